$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths (compensated for the COM pixel-rounding applied when
# translating the "characters" ColumnWidth into the stored OOXML width, so
# the persisted <col width="..."/> lands on the target value).
$ws.Columns.Item(1).ColumnWidth = 47.736979166666664
$ws.Columns.Item(2).ColumnWidth = 10.592447916666666
$ws.Columns.Item(3).ColumnWidth = 83.59244791666667
$ws.Columns.Item(5).ColumnWidth = 16.166666666666668

# New rows of data (rows 6-20)
$data = @(
    @("Two Sum", "Easy", "https://leetcode.com/problems/two-sum", "array", "UNSOLVED"),
    @("Median of Two Sorted Arrays", "Hard", "https://leetcode.com/problems/median-of-two-sorted-arrays", "array", "UNSOLVED"),
    @("Container With Most Water", "Medium", "https://leetcode.com/problems/container-with-most-water", "array", "UNSOLVED"),
    @("3Sum", "Medium", "https://leetcode.com/problems/3sum", "array", "UNSOLVED"),
    @("3Sum Closest", "Medium", "https://leetcode.com/problems/3sum-closest", "array", "UNSOLVED"),
    @("4Sum", "Medium", "https://leetcode.com/problems/4sum", "array", "UNSOLVED"),
    @("Remove Duplicates from Sorted Array", "Easy", "https://leetcode.com/problems/remove-duplicates-from-sorted-array", "array", "UNSOLVED"),
    @("Remove Element", "Easy", "https://leetcode.com/problems/remove-element", "array", "UNSOLVED"),
    @("Next Permutation", "Medium", "https://leetcode.com/problems/next-permutation", "array", "UNSOLVED"),
    @("Search Insert Position", "Easy", "https://leetcode.com/problems/search-insert-position", "array", "UNSOLVED"),
    @("Height Checker", "Easy", "https://leetcode.com/problems/height-checker", "counting-sort", "UNSOLVED"),
    @("Relative Sort Array", "Easy", "https://leetcode.com/problems/relative-sort-array", "counting-sort", "UNSOLVED"),
    @("Plus One", "Easy", "https://leetcode.com/problems/plus-one", "array", "UNSOLVED"),
    @("Pascal's Triangle II", "Easy", "https://leetcode.com/problems/pascals-triangle-ii", "array", "UNSOLVED"),
    @("Single Number", "Easy", "https://leetcode.com/problems/single-number", "array", "UNSOLVED")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Update selection to A15:E15 with active cell A15
$ws.Range("A15:E15").Select()
